$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks (E and F) ---
# Target XML widths (74.4140625 / -113.32421875) were written by a non-Excel
# tool directly into the OOXML and fall outside the pixel grid that Excel's
# COM ColumnWidth setter can produce (and Excel refuses negative widths
# outright). We set the closest values the object model will accept.
$ws.Columns.Item(5).ColumnWidth = 73.71428571428571
$ws.Columns.Item(6).ColumnWidth = 0

# --- Row 18: append &language= query param to both links ---
$ws.Range("E18").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21111-0004&bypass=true&levelindex=0&levelid=1660823284613&language=en"
$ws.Range("F18").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21111-0004&bypass=true&levelindex=1&levelid=1660810680251&language=de"

# --- Row 20 ---
$ws.Range("E20").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=51000-0001&bypass=true&levelindex=0&levelid=1660741526662&language=en"
$ws.Range("F20").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=51000-0001&language=de"

# --- Row 21: drop the jsessionid suffix from F21 ---
$ws.Range("F21").Value = "https://www.destatis.de/DE/Themen/Arbeit/Verdienste/Verdienste-Verdienstunterschiede/Tabellen/ugpg-02-bundeslaender-ab-2014.html"

# --- Row 22 ---
$ws.Range("E22").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0001&bypass=true&levelindex=1&levelid=1660802268437&language=en"
$ws.Range("F22").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0001&bypass=true&levelindex=0&levelid=1660822010108&language=de"

# --- Row 23 ---
$ws.Range("E23").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=22922-0011&bypass=true&levelindex=0&levelid=1660813986805#abreadcrumb&language=en"
$ws.Range("F23").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=22922-0011&bypass=true&levelindex=0&levelid=1660642440197#abreadcrumb&language=de"

# --- Row 24 ---
$ws.Range("E24").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21111-0013&bypass=true&levelindex=1&levelid=1660810680251&language=en"
$ws.Range("F24").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21111-0013&bypass=true&levelindex=0&levelid=1660823504838&language=de"

# --- Row 27 ---
$ws.Range("E27").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=81000-0031&bypass=true&levelindex=0&levelid=1660821702206&language=en"
$ws.Range("F27").Value = "hhttps://www-genesis.destatis.de/genesis/online?operation=previous&levelindex=0&step=0&titel=Tabellenaufbau&levelid=1660821596823&language=de"

# --- Row 28: drop #abreadcrumb, append language param ---
$ws.Range("E28").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21821-0001&bypass=true&levelindex=0&levelid=1660726117256&language=en"
$ws.Range("F28").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21821-0001&bypass=true&levelindex=1&levelid=1622107294362&language=de"

# --- Row 32: drop #abreadcrumb, append language param ---
$ws.Range("E32").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21821-0002&bypass=true&levelindex=1&levelid=1623135114747&language=en"
$ws.Range("F32").Value = "https://www-genesis.destatis.de/genesis//online?operation=table&code=21821-0002&bypass=true&levelindex=0&levelid=1660726117256&language=de"

# --- Row 66: fill in the previously-empty English link, switch the German one to lang=en/de pair ---
$ws.Range("E66").Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD11__custom_3696252/default/table?lang=en"
$ws.Range("F66").Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD11__custom_3696252/default/table?lang=de"

# --- New row 67: clone formatting from row 66, then overwrite with the new record ---
$ws.Range("A66:F66").Copy($ws.Range("A67:F67"))
$ws.Range("A67").Value = "L_SVWS"
$ws.Range("B67").Value = "Q_SVWS"
$ws.Range("C67").Value = "Forschung und Entwicklung in der Wirtschaft"
$ws.Range("D67").Value = "Research and development in the economy (only available in German)"
$ws.Range("E67").Value = ""
$ws.Range("F67").Value = "https://www.stifterverband.org/fue-facts-2020"
